# Update discount-factors sheet: label the maturities column with "NM"
# text (e.g. "1M", "2M", ...) instead of plain numbers, rename the header
# from "months" to "maturities", and tidy up the view state to match the
# re-uploaded workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: A1 "months" -> "maturities" (B1 "dfs" is left as-is).
$ws.Range("A1").Value = "maturities"

# Column A: same maturity numbers as before, now written as text with an
# "M" (months) suffix, e.g. 1 -> "1M", 720 -> "720M".
$maturities = @(1,2,3,4,5,6,7,8,9,10,11,12,15,18,21,24,36,48,60,72,84,96,108,120,132,144,180,240,300,360,480,600,720)

for ($i = 0; $i -lt $maturities.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = "$($maturities[$i])M"
}

# Column A width, close to the source file's best-fit width for the
# "maturities" header text.
$ws.Columns.Item(1).ColumnWidth = 9.3

# Selection state left behind by the last save.
$ws.Range("D1:D1048576").Select() | Out-Null

Write-Host "done"
